# Apply updated crypto price/volume figures (Mon Aug 21 04:31:27 UTC 2023 refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.277.64"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "1.688.73"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'217.41"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").Value = "'0.5349"
$ws.Range("E6").Value = "  +1.74%  "
$ws.Range("D8").Value = "'0.2713"
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("D9").Value = "'0.06412"
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("D10").Value = "'21.61"
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("D11").Value = "'0.07675"
$ws.Range("E11").Value = "  +2.17%  "
$ws.Range("D12").Value = "1.694.19"
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("D13").Value = "'4.527"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").Value = "'0.5772"
$ws.Range("E14").Value = "  -0.59%  "
$ws.Range("D15").Value = "'0.000008352"
$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("D16").Value = "'66.49"
$ws.Range("E16").Value = "  +2.65%  "
$ws.Range("D17").Value = "26.307.73"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D19").Value = "'4.880"
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("D21").Value = "'190.72"
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("D22").Value = "'6.245"
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'148.54"
$ws.Range("E24").Value = "  +2.40%  "
$ws.Range("D25").Value = "'0.1286"
$ws.Range("E25").Value = "  +2.75%  "
$ws.Range("D26").Value = "'7.845"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").Value = "'15.79"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").Value = "'0.06172"
$ws.Range("E28").Value = "  -4.71%  "
$ws.Range("D29").Value = "'1.374"
$ws.Range("E29").Value = "  +0.92%  "
$ws.Range("D30").Value = "'1.324"
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("D31").Value = "'3.593"
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("D32").Value = "'3.579"
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("E33").Value = "  +0.51%  "
$ws.Range("E34").Value = "  -0.23%  "
$ws.Range("D35").Value = "'0.6198"
$ws.Range("E35").Value = "  -0.56%  "
$ws.Range("D36").Value = "'2.424"
$ws.Range("E36").Value = "  +0.67%  "
$ws.Range("E37").Value = "  +0.71%  "
$ws.Range("D38").Value = "'0.01647"
$ws.Range("E38").Value = "  +1.56%  "
$ws.Range("E39").Value = "  -4.66%  "
$ws.Range("D40").Value = "1.106.92"
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("D41").Value = "'0.8807"
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("D42").Value = "'1.013"
$ws.Range("E42").Value = "  -0.26%  "
$ws.Range("D43").Value = "'100.82"
$ws.Range("E43").Value = "  +0.34%  "
$ws.Range("D44").Value = "1.839.52"
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("E45").Value = "  +2.13%  "
$ws.Range("D46").Value = "'57.54"
$ws.Range("E46").Value = "  +1.10%  "
$ws.Range("D47").Value = "'8.159"
$ws.Range("E47").Value = "  -0.91%  "
$ws.Range("D48").Value = "'1.002"
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("D49").Value = "'0.05283"
$ws.Range("D50").Value = "'0.4300"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").Value = "'6.052"
$ws.Range("E51").Value = "  -0.42%  "
